$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the Coffee_Type labels in column B were swapped (DECAF rows were
# actually CAF and vice versa). Correct each block of rows to the right label.

# Rows that were mislabeled "DECAF" -> should be "CAF"
$ws.Range("B2:B17").Value = "CAF"
$ws.Range("B33:B48").Value = "CAF"
$ws.Range("B64:B79").Value = "CAF"
$ws.Range("B95:B110").Value = "CAF"

# Rows that were mislabeled "CAF" -> should be "DECAF"
$ws.Range("B18:B32").Value = "DECAF"
$ws.Range("B49:B63").Value = "DECAF"
$ws.Range("B80:B94").Value = "DECAF"
$ws.Range("B111:B125").Value = "DECAF"

# Restore the selection left on the sheet by the editor (cell C73)
$ws.Range("C73").Select()
